$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "resume_test - Copy - Copy - Copy.xyx"
$ws.Range("B2").Value = "unsupported file (not a .pdf nor .docx"
$ws.Range("C2").Value = "unsupported file (not a .pdf nor .docx"
$ws.Range("D2").Value = "unsupported file (not a .pdf nor .docx"
